$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030122438832786
$ws.Range("D2").Value = 1.033003538817872
$ws.Range("E2").Value = 1.038556920987629
$ws.Range("F2").Value = 1.046434912917202
$ws.Range("I2").Value = 1.030617215237833
$ws.Range("J2").Value = 1.035265428256908
$ws.Range("K2").Value = 1.03580713958557
$ws.Range("L2").Value = 1.041344615790479
$ws.Range("M2").Value = 1.049200357856858
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031149140743714
$ws.Range("D3").Value = 1.03375066270984
$ws.Range("E3").Value = 1.039505178173775
$ws.Range("F3").Value = 1.047551122915422
$ws.Range("I3").Value = 1.030760051608249
$ws.Range("J3").Value = 1.035933056745677
$ws.Range("K3").Value = 1.036363318453543
$ws.Range("L3").Value = 1.042102545651459
$ws.Range("M3").Value = 1.050127418465016
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031813752364017
$ws.Range("D4").Value = 1.034234010266911
$ws.Range("E4").Value = 1.040119421221849
$ws.Range("F4").Value = 1.048274367238374
$ws.Range("I4").Value = 1.030851057286696
$ws.Range("J4").Value = 1.036364755063837
$ws.Range("K4").Value = 1.036722449190801
$ws.Range("L4").Value = 1.04259301276967
$ws.Range("M4").Value = 1.050727682543994
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032093218738324
$ws.Range("D5").Value = 1.034437186662235
$ws.Range("E5").Value = 1.040377805673
$ws.Range("F5").Value = 1.048578653724663
$ws.Range("I5").Value = 1.030888975958084
$ws.Range("J5").Value = 1.036546168222429
$ws.Range("K5").Value = 1.036873246501165
$ws.Range("L5").Value = 1.042799213188205
$ws.Range("M5").Value = 1.050980127361645
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032140146098097
$ws.Range("D6").Value = 1.034471299500084
$ws.Range("E6").Value = 1.040421198727258
$ws.Range("F6").Value = 1.04862975850637
$ws.Range("I6").Value = 1.030895322713586
$ws.Range("J6").Value = 1.036576624004903
$ws.Range("K6").Value = 1.036898555384107
$ws.Range("L6").Value = 1.04283383560195
$ws.Range("M6").Value = 1.051022519443412
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031817486357571
$ws.Range("D7").Value = 1.03423672521116
$ws.Range("E7").Value = 1.04012287315275
$ws.Range("F7").Value = 1.048278432208176
$ws.Range("I7").Value = 1.030851565294767
$ws.Range("J7").Value = 1.036367179402111
$ws.Range("K7").Value = 1.036724464864919
$ws.Range("L7").Value = 1.042595767998822
$ws.Range("M7").Value = 1.050731055356055
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030469362162039
$ws.Range("D8").Value = 1.033256051193803
$ws.Range("E8").Value = 1.038877252453745
$ws.Range("F8").Value = 1.046811937734016
$ws.Range("I8").Value = 1.030665781029764
$ws.Range("J8").Value = 1.035491118636849
$ws.Range("K8").Value = 1.035995258576101
$ws.Range("L8").Value = 1.041600753999041
$ws.Range("M8").Value = 1.049513580476975
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028095848571114
$ws.Range("D9").Value = 1.031527320859324
$ws.Range("E9").Value = 1.036687372711874
$ws.Range("F9").Value = 1.04423532739616
$ws.Range("I9").Value = 1.030327555264267
$ws.Range("J9").Value = 1.033945097979429
$ws.Range("K9").Value = 1.034704558457858
$ws.Range("L9").Value = 1.039847709820645
$ws.Range("M9").Value = 1.047371272453085
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026514897574516
$ws.Range("D10").Value = 1.030374449674563
$ws.Range("E10").Value = 1.03523089481812
$ws.Range("F10").Value = 1.042522682524821
$ws.Range("I10").Value = 1.030094799734695
$ws.Range("J10").Value = 1.032912907518877
$ws.Range("K10").Value = 1.033840268825286
$ws.Range("L10").Value = 1.038679243567145
$ws.Range("M10").Value = 1.045945137061213
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025830659779186
$ws.Range("D11").Value = 1.029875164304372
$ws.Range("E11").Value = 1.03460104756962
$ws.Range("F11").Value = 1.041782301603727
$ws.Range("I11").Value = 1.029992294000037
$ws.Range("J11").Value = 1.032465605186363
$ws.Range("K11").Value = 1.033465123151509
$ws.Range("L11").Value = 1.038173345033973
$ws.Range("M11").Value = 1.0453281007404
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025576552146557
$ws.Range("D12").Value = 1.02968969567695
$ws.Range("E12").Value = 1.034367217572923
$ws.Range("F12").Value = 1.041507472584186
$ws.Range("I12").Value = 1.029953960578202
$ws.Range("J12").Value = 1.032299403932162
$ws.Range("K12").Value = 1.033325642331046
$ws.Range("L12").Value = 1.037985440460702
$ws.Range("M12").Value = 1.045098979882269
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02563105689009
$ws.Range("D13").Value = 1.029729479847276
$ws.Range("E13").Value = 1.034417369312221
$ws.Range("F13").Value = 1.041566416127466
$ws.Range("I13").Value = 1.02996219490719
$ws.Range("J13").Value = 1.032335057068835
$ws.Range("K13").Value = 1.033355567540779
$ws.Range("L13").Value = 1.038025746230979
$ws.Range("M13").Value = 1.045148123734414
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025809654168438
$ws.Range("D14").Value = 1.029859833645552
$ws.Range("E14").Value = 1.034581716590502
$ws.Range("F14").Value = 1.041759580448629
$ws.Range("I14").Value = 1.029989130613227
$ws.Range("J14").Value = 1.032451868022896
$ws.Range("K14").Value = 1.033453596374235
$ws.Range("L14").Value = 1.038157812604237
$ws.Range("M14").Value = 1.045309160029229
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025919700306341
$ws.Range("D15").Value = 1.029940147384532
$ws.Range("E15").Value = 1.034682992729057
$ws.Range("F15").Value = 1.041878619410801
$ws.Range("I15").Value = 1.030005692383213
$ws.Range("J15").Value = 1.03252383206229
$ws.Range("K15").Value = 1.03351397730075
$ws.Range("L15").Value = 1.038239184215834
$ws.Range("M15").Value = 1.045408389608405
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026560315019395
$ws.Range("D16").Value = 1.030407583895161
$ws.Range("E16").Value = 1.03527271296995
$ws.Range("F16").Value = 1.042571844593721
$ws.Range("I16").Value = 1.030101566471174
$ws.Range("J16").Value = 1.032942586012231
$ws.Range("K16").Value = 1.033865147042236
$ws.Range("L16").Value = 1.038712819580425
$ws.Range("M16").Value = 1.045986098075572
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026962242566145
$ws.Range("D17").Value = 1.030700772464746
$ws.Range("E17").Value = 1.035642848080297
$ws.Range("F17").Value = 1.043007009416558
$ws.Range("I17").Value = 1.030161245219143
$ws.Range("J17").Value = 1.033205163984797
$ws.Range("K17").Value = 1.034085185316039
$ws.Range("L17").Value = 1.039009933551801
$ws.Range("M17").Value = 1.046348610677633
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027196711404921
$ws.Range("D18").Value = 1.030871776180561
$ws.Range("E18").Value = 1.0358588204493
$ws.Range("F18").Value = 1.043260949999512
$ws.Range("I18").Value = 1.030195888743459
$ws.Range("J18").Value = 1.033358286797126
$ws.Range("K18").Value = 1.034213442795246
$ws.Range("L18").Value = 1.039183240391026
$ws.Range("M18").Value = 1.046560105373249
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027276664528646
$ws.Range("D19").Value = 1.030930082607738
$ws.Range("E19").Value = 1.035932474799364
$ws.Range("F19").Value = 1.043347556894301
$ws.Range("I19").Value = 1.030207673109617
$ws.Range("J19").Value = 1.033410491859368
$ws.Range("K19").Value = 1.034257160485773
$ws.Range("L19").Value = 1.039242334429264
$ws.Range("M19").Value = 1.046632227599615
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026919116306414
$ws.Range("D20").Value = 1.030669316949313
$ws.Range("E20").Value = 1.035603127929144
$ws.Range("F20").Value = 1.042960308317344
$ws.Range("I20").Value = 1.030154859430195
$ws.Range("J20").Value = 1.033176995416539
$ws.Range("K20").Value = 1.034061586298952
$ws.Range("L20").Value = 1.038978055484809
$ws.Range("M20").Value = 1.046309711618725
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025757060375539
$ws.Range("D21").Value = 1.029821448019818
$ws.Range("E21").Value = 1.034533317004954
$ws.Range("F21").Value = 1.041702693360109
$ws.Range("I21").Value = 1.029981205844767
$ws.Range("J21").Value = 1.032417471575225
$ws.Range("K21").Value = 1.03342473304843
$ws.Range("L21").Value = 1.038118922106554
$ws.Range("M21").Value = 1.045261736811322
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025026711751945
$ws.Range("D22").Value = 1.029288291281232
$ws.Range("E22").Value = 1.033861398142577
$ws.Range("F22").Value = 1.040913030708698
$ws.Range("I22").Value = 1.029870529085378
$ws.Range("J22").Value = 1.031939620886945
$ws.Range("K22").Value = 1.03302353727427
$ws.Range("L22").Value = 1.037578801194078
$ws.Range("M22").Value = 1.044603261213015
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02541385663534
$ws.Range("D23").Value = 1.02957093389433
$ws.Range("E23").Value = 1.034217527215449
$ws.Range("F23").Value = 1.041331546215502
$ws.Range("I23").Value = 1.029929342416101
$ws.Range("J23").Value = 1.032192967632886
$ws.Range("K23").Value = 1.033236292566044
$ws.Range("L23").Value = 1.037865124605453
$ws.Range("M23").Value = 1.04495229092888
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026938603109988
$ws.Range("D24").Value = 1.030683530368174
$ws.Range("E24").Value = 1.035621075512142
$ws.Range("F24").Value = 1.042981410174409
$ws.Range("I24").Value = 1.030157745406514
$ws.Range("J24").Value = 1.033189723686521
$ws.Range("K24").Value = 1.034072249948158
$ws.Range("L24").Value = 1.038992459794731
$ws.Range("M24").Value = 1.046327288283061
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028709214283404
$ws.Range("D25").Value = 1.031974311040476
$ws.Range("E25").Value = 1.037252905223277
$ws.Range("F25").Value = 1.044900546729857
$ws.Range("I25").Value = 1.03041627812603
$ws.Range("J25").Value = 1.034345049684559
$ws.Range("K25").Value = 1.035038912115195
$ws.Range("L25").Value = 1.040300875428463
$ws.Range("M25").Value = 1.047924747583108

Write-Host "Updated $(241) cells"
